# Auto-generated Excel COM-interop edit script
# Applies the cell-value changes described by the authoritative XML diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (67 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 271.8889
$ws.Range("J9").Value = 174.66667
$ws.Range("L9").Value = 174.66667
$ws.Range("N9").Value = -512.6666700000001
$ws.Range("H13").Value = 110981
$ws.Range("I13").Value = 367201.34
$ws.Range("K13").Value = 367201.34
$ws.Range("M13").Value = -367032.34
$ws.Range("H15").Value = 813.7692
$ws.Range("I15").Value = 813.7692
$ws.Range("K15").Value = 2441.3076
$ws.Range("M15").Value = -2272.3076
$ws.Range("H17").Value = 436631.44
$ws.Range("J17").Value = 436631.44
$ws.Range("L17").Value = 1309894.32
$ws.Range("N17").Value = -1310230.32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H107").Value = 96326.625
$ws.Range("I107").Value = 1323
$ws.Range("J107").Value = 254666
$ws.Range("K107").Value = 1323
$ws.Range("L107").Value = 254666
$ws.Range("M107").Value = 597
$ws.Range("N107").Value = -258506
$ws.Range("H116").Value = 4569.4287
$ws.Range("I116").Value = 3993
$ws.Range("J116").Value = 4800
$ws.Range("K116").Value = 3993
$ws.Range("L116").Value = 4800
$ws.Range("M116").Value = -551
$ws.Range("N116").Value = -11684
$ws.Range("H125").Value = 1830.1428
$ws.Range("I125").Value = 1453
$ws.Range("K125").Value = 13077
$ws.Range("M125").Value = -10617
$ws.Range("H127").Value = 667.4286
$ws.Range("I127").Value = 667.4286
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2002.2858
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2957.7142
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 1637.85
$ws.Range("I129").Value = 717.4666999999999
$ws.Range("K129").Value = 2152.4001
$ws.Range("M129").Value = 2847.5999
$ws.Range("H132").Value = 1881.68
$ws.Range("I132").Value = 1881.68
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5645.04
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3115.04
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 2586.1667
$ws.Range("I138").Value = 3774.2
$ws.Range("J138").Value = 2129.2307
$ws.Range("K138").Value = 11322.6
$ws.Range("L138").Value = 6387.6921
$ws.Range("M138").Value = -6182.599999999999
$ws.Range("N138").Value = -16667.6921
$ws.Range("H141").Value = 7875.9473
$ws.Range("I141").Value = 7887.9165
$ws.Range("K141").Value = 23663.7495
$ws.Range("M141").Value = -18483.7495

# ---- Sheet: ARM (57 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 63.666668
$ws.Range("I5").Value = 63.666668
$ws.Range("K5").Value = 63.666668
$ws.Range("M5").Value = 48.333332
$ws.Range("H8").Value = 12799.833
$ws.Range("I8").Value = 8499.5
$ws.Range("J8").Value = 14950
$ws.Range("K8").Value = 8499.5
$ws.Range("L8").Value = 14950
$ws.Range("M8").Value = -8355.5
$ws.Range("N8").Value = -15238
$ws.Range("H10").Value = 4000
$ws.Range("I10").Value = 2750
$ws.Range("J10").Value = 4833.3335
$ws.Range("K10").Value = 2750
$ws.Range("L10").Value = 4833.3335
$ws.Range("M10").Value = -2580
$ws.Range("N10").Value = -5173.3335
$ws.Range("H11").Value = 8669.333000000001
$ws.Range("J11").Value = 8669.333000000001
$ws.Range("L11").Value = 8669.333000000001
$ws.Range("N11").Value = -8957.333000000001
$ws.Range("H25").Value = 3031.6
$ws.Range("I25").Value = 857
$ws.Range("K25").Value = 857
$ws.Range("M25").Value = -455
$ws.Range("H61").Value = 2498.2
$ws.Range("I61").Value = 2193.3914
$ws.Range("K61").Value = 2193.3914
$ws.Range("M61").Value = -1981.3914
$ws.Range("H74").Value = 1517.9464
$ws.Range("I74").Value = 1336.42
$ws.Range("K74").Value = 1336.42
$ws.Range("M74").Value = -462.4200000000001
$ws.Range("H77").Value = 1517.9464
$ws.Range("I77").Value = 1336.42
$ws.Range("K77").Value = 6682.1
$ws.Range("M77").Value = -2314.1
$ws.Range("H102").Value = 4374
$ws.Range("I102").Value = 2056.25
$ws.Range("K102").Value = 2056.25
$ws.Range("M102").Value = -434.25
$ws.Range("H122").Value = 2993.25
$ws.Range("I122").Value = 3056.7334
$ws.Range("J122").Value = 2802.8
$ws.Range("K122").Value = 9170.200199999999
$ws.Range("L122").Value = 8408.400000000001
$ws.Range("M122").Value = -6720.200199999999
$ws.Range("N122").Value = -13308.4
$ws.Range("H136").Value = 2498.2
$ws.Range("I136").Value = 2193.3914
$ws.Range("K136").Value = 6580.174199999999
$ws.Range("M136").Value = -4030.174199999999

# ---- Sheet: BSM (45 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 63.666668
$ws.Range("I4").Value = 63.666668
$ws.Range("K4").Value = 63.666668
$ws.Range("M4").Value = 51.333332
$ws.Range("H20").Value = 15231.6
$ws.Range("I20").Value = 14073.5
$ws.Range("J20").Value = 17547.8
$ws.Range("K20").Value = 14073.5
$ws.Range("L20").Value = 17547.8
$ws.Range("M20").Value = -13826.5
$ws.Range("N20").Value = -18041.8
$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 2000
$ws.Range("N39").Value = -2778
$ws.Range("H86").Value = 1952.1111
$ws.Range("I86").Value = 1652.6
$ws.Range("J86").Value = 2326.5
$ws.Range("K86").Value = 1652.6
$ws.Range("L86").Value = 2326.5
$ws.Range("M86").Value = -529.5999999999999
$ws.Range("N86").Value = -4572.5
$ws.Range("H89").Value = 1952.1111
$ws.Range("I89").Value = 1652.6
$ws.Range("J89").Value = 2326.5
$ws.Range("K89").Value = 8263
$ws.Range("L89").Value = 11632.5
$ws.Range("M89").Value = -2647
$ws.Range("N89").Value = -22864.5
$ws.Range("H94").Value = 3389.1785
$ws.Range("I94").Value = 3173.8262
$ws.Range("K94").Value = 3173.8262
$ws.Range("M94").Value = -2722.8262
$ws.Range("H95").Value = 22251.666
$ws.Range("J95").Value = 22251.666
$ws.Range("L95").Value = 22251.666
$ws.Range("N95").Value = -27743.666
$ws.Range("H107").Value = 60773.59
$ws.Range("I107").Value = 92778.09
$ws.Range("K107").Value = 92778.09
$ws.Range("M107").Value = -90858.09
$ws.Range("H134").Value = 925.15
$ws.Range("I134").Value = 806.05884
$ws.Range("K134").Value = 2418.17652
$ws.Range("M134").Value = 116.82348

# ---- Sheet: CRP (65 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 500001000
$ws.Range("I7").Value = 500001000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 500001000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -500000887
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 1258.7
$ws.Range("I22").Value = 1248.6
$ws.Range("K22").Value = 1248.6
$ws.Range("M22").Value = -898.5999999999999
$ws.Range("H42").Value = 7833.3335
$ws.Range("I42").Value = 1000
$ws.Range("J42").Value = 11250
$ws.Range("K42").Value = 1000
$ws.Range("L42").Value = 11250
$ws.Range("M42").Value = -407
$ws.Range("N42").Value = -12436
$ws.Range("H58").Value = 300
$ws.Range("I58").Value = 300
$ws.Range("K58").Value = 300
$ws.Range("M58").Value = -97
$ws.Range("H86").Value = 4992
$ws.Range("I86").Value = 4989
$ws.Range("K86").Value = 4989
$ws.Range("M86").Value = -3866
$ws.Range("H89").Value = 4992
$ws.Range("I89").Value = 4989
$ws.Range("K89").Value = 24945
$ws.Range("M89").Value = -19329
$ws.Range("H96").Value = 13084.4
$ws.Range("J96").Value = 13084.4
$ws.Range("L96").Value = 13084.4
$ws.Range("N96").Value = -18576.4
$ws.Range("H99").Value = 65555.336
$ws.Range("I99").Value = 65555.336
$ws.Range("K99").Value = 65555.336
$ws.Range("M99").Value = -64057.336
$ws.Range("H107").Value = 1341.3243
$ws.Range("I107").Value = 1058.4286
$ws.Range("K107").Value = 1058.4286
$ws.Range("M107").Value = 861.5714
$ws.Range("H122").Value = 255449.75
$ws.Range("I122").Value = 255449.75
$ws.Range("K122").Value = 766349.25
$ws.Range("M122").Value = -763899.25
$ws.Range("H126").Value = 65555.336
$ws.Range("I126").Value = 65555.336
$ws.Range("K126").Value = 196666.008
$ws.Range("M126").Value = -194196.008
$ws.Range("H132").Value = 2908.7104
$ws.Range("I132").Value = 3014.08
$ws.Range("K132").Value = 9042.24
$ws.Range("M132").Value = -6512.24
$ws.Range("H134").Value = 3104.5715
$ws.Range("I134").Value = 3257.1667
$ws.Range("J134").Value = 2189
$ws.Range("K134").Value = 9771.500100000001
$ws.Range("L134").Value = 6567
$ws.Range("M134").Value = -7236.500100000001
$ws.Range("N134").Value = -11637
$ws.Range("H136").Value = 300
$ws.Range("I136").Value = 300
$ws.Range("K136").Value = 900
$ws.Range("M136").Value = 1650

# ---- Sheet: CUL (48 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14671468
$ws.Range("I4").Value = 2273154.8
$ws.Range("J4").Value = 37401710
$ws.Range("K4").Value = 6819464.399999999
$ws.Range("L4").Value = 112205130
$ws.Range("M4").Value = -6819352.399999999
$ws.Range("N4").Value = -112205354
$ws.Range("H5").Value = 470.875
$ws.Range("J5").Value = 305
$ws.Range("L5").Value = 915
$ws.Range("N5").Value = -1139
$ws.Range("H11").Value = 1222651.9
$ws.Range("I11").Value = 144.83333
$ws.Range("J11").Value = 3667666
$ws.Range("K11").Value = 434.49999
$ws.Range("L11").Value = 11002998
$ws.Range("M11").Value = -294.49999
$ws.Range("N11").Value = -11003278
$ws.Range("H34").Value = 281.83334
$ws.Range("J34").Value = 995
$ws.Range("L34").Value = 2985
$ws.Range("N34").Value = -3153
$ws.Range("H39").Value = 3560
$ws.Range("J39").Value = 3560
$ws.Range("L39").Value = 10680
$ws.Range("N39").Value = -11268
$ws.Range("H55").Value = 7354230
$ws.Range("J55").Value = 13890189
$ws.Range("L55").Value = 41670567
$ws.Range("N55").Value = -41670921
$ws.Range("H103").Value = 466.42856
$ws.Range("I103").Value = 477.5
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 1432.5
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = -553.5
$ws.Range("N103").Value = -2958
$ws.Range("H135").Value = 470.875
$ws.Range("J135").Value = 305
$ws.Range("L135").Value = 2745
$ws.Range("N135").Value = -7815
$ws.Range("H140").Value = 3270.923
$ws.Range("I140").Value = 1399.5
$ws.Range("J140").Value = 6265.2
$ws.Range("K140").Value = 4198.5
$ws.Range("L140").Value = 18795.6
$ws.Range("M140").Value = 981.5
$ws.Range("N140").Value = -29155.6

# ---- Sheet: GSM (19 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4000689.8
$ws.Range("J3").Value = 5001500
$ws.Range("L3").Value = 5001500
$ws.Range("N3").Value = -5001732
$ws.Range("H102").Value = 2825.2
$ws.Range("I102").Value = 2697
$ws.Range("K102").Value = 2697
$ws.Range("M102").Value = -1075
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 1996.7693
$ws.Range("I122").Value = 1913.4166
$ws.Range("J122").Value = 2997
$ws.Range("K122").Value = 5740.2498
$ws.Range("L122").Value = 8991
$ws.Range("M122").Value = -3290.2498
$ws.Range("N122").Value = -13891

# ---- Sheet: LTW (31 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11723.433
$ws.Range("I7").Value = 21875.857
$ws.Range("J7").Value = 5543.696
$ws.Range("K7").Value = 21875.857
$ws.Range("L7").Value = 5543.696
$ws.Range("M7").Value = -21763.857
$ws.Range("N7").Value = -5767.696
$ws.Range("H34").Value = 2999.6667
$ws.Range("H40").Value = 6522.3335
$ws.Range("I40").Value = 4886.5557
$ws.Range("K40").Value = 4886.5557
$ws.Range("M40").Value = -4750.5557
$ws.Range("H43").Value = 69696970
$ws.Range("J43").Value = 69696970
$ws.Range("L43").Value = 69696970
$ws.Range("N43").Value = -69697356
$ws.Range("H93").Value = 34416.91
$ws.Range("I93").Value = 4181.125
$ws.Range("K93").Value = 4181.125
$ws.Range("M93").Value = -2933.125
$ws.Range("H126").Value = 11723.433
$ws.Range("I126").Value = 21875.857
$ws.Range("J126").Value = 5543.696
$ws.Range("K126").Value = 65627.571
$ws.Range("L126").Value = 16631.088
$ws.Range("M126").Value = -63157.571
$ws.Range("N126").Value = -21571.088
$ws.Range("H136").Value = 4780.3687
$ws.Range("I136").Value = 4217.923
$ws.Range("K136").Value = 12653.769
$ws.Range("M136").Value = -10103.769

# ---- Sheet: WVR (15 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 9250.833000000001
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 11001
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 11001
$ws.Range("M7").Value = -387
$ws.Range("N7").Value = -11227
$ws.Range("H132").Value = 4610.2974
$ws.Range("I132").Value = 5059.483
$ws.Range("K132").Value = 15178.449
$ws.Range("M132").Value = -12648.449
$ws.Range("H136").Value = 1311
$ws.Range("I136").Value = 1009.5455
$ws.Range("K136").Value = 3028.6365
$ws.Range("M136").Value = -478.6364999999996
